$wb = $excel.ActiveWorkbook

function Update-Sheet($ws, $BValues, $CValues) {
    for ($i = 0; $i -lt $BValues.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = ($i + 1)
        $ws.Cells.Item($row, 2).Value = $BValues[$i]
        $ws.Cells.Item($row, 3).Value = $CValues[$i]
    }
}

# Sheet2 - sigma_010
$ws2 = $wb.Worksheets.Item("sigma_010")
Update-Sheet $ws2 `
    @(28.20305651627879,28.1580144094411,28.16987348421478,28.15600956894925,28.16395428689887,28.16339661077645,28.14633007291994,28.20476797677451,28.18092803791902,28.17665405271457) `
    @(33.00751797920881,32.99401465644112,32.99175526194325,33.02847833869218,32.97946245982602,33.00158832012826,32.98248055393285,33.01067652988006,33.00732982539505,32.99690069897951)
$ws2.Cells.Item(12, 2).Value = 28.17229850168873
$ws2.Cells.Item(12, 3).Value = 33.00002046244271

# Sheet3 - sigma_025
$ws3 = $wb.Worksheets.Item("sigma_025")
Update-Sheet $ws3 `
    @(19.66903947027433,19.66565577383501,19.67155795443999,19.66982823375871,19.6632518252504,19.65307407468944,19.67157095610282,19.66022463844488,19.65748450088456,19.66824831209693) `
    @(30.05712060327577,30.06586323507985,30.03929326543123,30.01422176071177,30.04131956310932,30.02211282573864,30.07542490954676,30.04663974566906,30.08895833902996,30.04012679888558)
$ws3.Cells.Item(12, 2).Value = 19.66499357397771
$ws3.Cells.Item(12, 3).Value = 30.04910810464779

# Sheet4 - sigma_050
$ws4 = $wb.Worksheets.Item("sigma_050")
Update-Sheet $ws4 `
    @(14.47884779077478,14.4995165635431,14.48510080738422,14.48497982696361,14.48706825028646,14.51690813240413,14.50663914390977,14.48977383929543,14.5059729093162,14.4933240800243) `
    @(25.12171099451316,25.11024588457892,25.11170989818008,25.14808107345671,25.06867512805906,25.12824280903553,25.14130665025453,25.0607769853883,25.07467567094465,25.04972300678268)
$ws4.Cells.Item(12, 2).Value = 14.4948131343902
$ws4.Cells.Item(12, 3).Value = 25.10151481011936
